# Add 2022-Q3 data: a new sheet "2022-Q3" inserted right after "2022-Q2" column's
# template, positioned before "2022-Q2" in the tab order; and the summary sheet
# "总计" gets a new leading row for 2022-Q3 plus the newly-revealed 2021-Q1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Build the new "2022-Q3" worksheet.
#    Use the existing "2022-Q2" sheet as a formatting template: copy it (so
#    header/index-column styles match exactly), place the copy immediately
#    before "2022-Q2", rename it, then overwrite the copied cell values with
#    the 2022-Q3 figures. A 13th row is added (the Q3 sheet has 12 data rows
#    vs. 11 on the Q2 sheet) by copying row 12's formatting down to row 13.
# ---------------------------------------------------------------------------

$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Extend formatting from the last existing row (12) down to the new row (13).
$q3.Range("A12:H12").Copy()
$q3.Range("A13:H13").PasteSpecial(-4122)

# Columns B:G hold values that look numeric ("011429", "60.79", ...) but must
# stay text (leading zeros / exact decimal text), matching the source data.
$q3.Range("B2:G13").NumberFormat = "@"

$q3Data = @(
    @(0,  "011429", "前海开源民裕进取混合",     "2.53", "60.79", "4.73", "0.1197", 2),
    @(1,  "630015", "华商大盘量化精选混合",     "2.68", "91.58", "2.75", "0.0737", 6),
    @(2,  "008488", "华商恒益稳健混合",         "2.39", "52.53", "1.97", "0.0471", 7),
    @(3,  "011588", "前海开源成份精选混合",     "0.84", "60.90", "4.47", "0.0375", 2),
    @(4,  "007251", "广发睿享稳健增利混合A",    "1.45", "32.79", "1.64", "0.0238", 10),
    @(5,  "001115", "广发聚安混合A",           "2.02", "21.12", "0.65", "0.0131", 7),
    @(6,  "006890", "上投摩根领先优选混合",     "0.32", "80.46", "3.80", "0.0122", 2),
    @(7,  "001252", "中海进取收益灵活配置混合", "0.20", "87.95", "4.11", "0.0082", 10),
    @(8,  "001116", "广发聚安混合C",           "1.17", "21.12", "0.65", "0.0076", 7),
    @(9,  "011702", "广发睿享稳健增利混合C",    "0.29", "32.79", "1.64", "0.0048", 10),
    @(10, "001914", "中信建投聚利混合A",        "0.10", "39.73", "2.05", "0.0020", 7),
    @(11, "006845", "中信建投聚利混合C",        "0.01", "39.73", "2.05", "0.0002", 7)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Drop the forced text number format again (the cells stay text because they
# were entered while formatted as text; clearing the format afterwards keeps
# that quoting but removes the now-redundant explicit style).
$q3.Range("B2:G13").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q3 totals as the new
#    row 2 and push the rest down, revealing a 2021-Q1 row at the bottom.
# ---------------------------------------------------------------------------

$summary = $wb.Worksheets.Item("总计")

$summary.Range("A2").Copy()
$summary.Range("A6").PasteSpecial(-4122)

$summaryRows = @(
    @("2022-Q3", 12, 0.35),
    @("2022-Q2", 11, 0.57),
    @("2022-Q1", 4,  0.17),
    @("2021-Q3", 2,  0.11),
    @("2021-Q1", 2,  0.01)
)

$r = 2
$idx = 0
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $idx
    $summary.Cells.Item($r, 2).Value = $row[0]
    $summary.Cells.Item($r, 3).Value = $row[1]
    $summary.Cells.Item($r, 4).Value = $row[2]
    $r++
    $idx++
}

Write-Output "done"
